# Corrected bug when counting number of completed flights.
# (Especially affecting D2C2 scenarios.)
#
# This updates the raw counts/sums on the "Low traffic densit" sheet
# (row 17, the first data row), and lets Excel recalculate the
# dependent formulas (F17 ratio, and the AVERAGE/STDEV.P summary rows).
# It also restores the active-sheet / selection state to match.

$wb = $excel.ActiveWorkbook

$wsLow = $wb.Worksheets.Item("Low traffic densit")
$wsMedium = $wb.Worksheets.Item("Medium traffic density")

# Corrected source data for row 17.
$wsLow.Range("C17").Value = 21
$wsLow.Range("D17").Value = 4443.81
$wsLow.Range("E17").Value = 458.45

# Recalculate all formulas (F17 ratio, AVERAGE/STDEV.P summary rows, etc.)
$excel.Calculate()

# Restore selection / active sheet state.
$wsMedium.Activate()
$wsMedium.Range("C24").Select() | Out-Null
$wsMedium.Tab.Selected = $false

$wsLow.Activate()
$wsLow.Range("C18").Select() | Out-Null
$wsLow.Tab.Selected = $true
